$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 41, shifting existing rows 41-75 down to 42-76
$ws.Rows("41:41").Insert()

# Populate the new row 41 with this week's data, copying the constant
# columns (A,B,C,E,F,G,H,I,N,O,Q,R) from the row below (now row 42, the
# former row 41) and setting the new weekly values (D,J,K,L,M,P).
$ws.Range("A41").Value = 8
$ws.Range("B41").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C41").Value = 'Coquimbo'
$ws.Range("D41").NumberFormat = $ws.Range("D42").NumberFormat
$ws.Range("D41").Value = "2022-11-09"
$ws.Range("E41").Value = 4
$ws.Range("F41").Value = 100114007
$ws.Range("G41").Value = 'Jengibre'
$ws.Range("H41").Value = 'Sin especificar'
$ws.Range("I41").Value = 'Primera'
$ws.Range("J41").Value = 520
$ws.Range("K41").Value = 13500
$ws.Range("L41").Value = 14000
$ws.Range("M41").Value = 13750
$ws.Range("N41").Value = '$/caja 13 kilos'
$ws.Range("O41").Value = 'Perú'
$ws.Range("P41").Value = 1058
$ws.Range("Q41").Value = 13
$ws.Range("R41").Value = 'Hortaliza'
